$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.643.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.84%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -7.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.92%  "

# Row 5
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.79%  "

# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.88"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3621"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.87%  "

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.70"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -8.88%  "

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3274"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -10.30%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -9.83%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07018"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -8.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.000"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.85%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.45"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -11.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.628"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.644.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.88%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001056"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -9.82%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06566"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.36%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.55"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -10.44%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -9.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.969"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -9.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.92%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.662.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.59%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.58%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.496"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -16.93%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.82"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "126.92"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.825.31"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.079"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -10.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.108"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.929"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -19.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.719"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.44%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08424"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.68%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.37"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -12.20%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.179"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06156"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -9.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02261"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -9.56%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.262"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -12.46%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2064"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.09%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.209"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.28%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.63%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5980"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -9.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.741"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.82%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.99"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5695"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -10.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.73"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -8.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.963"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -10.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07048"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.40"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -8.33%  "
